$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E) to (B:F)
$ws.Columns("A:A").Insert()

# Populate the new ID column (header + row labels)
$ids = @(
    "ID",
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

# Apply the same header style (bold/bordered/centered) used by the other headers to the new ID header cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
